$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: E14 was stored as text "+919328027733"; it becomes a plain number.
$ws.Range("E14").Value = 919328027733

# Row 15 (new)
$ws.Range("A15").Value = "Yes, I am there."
$ws.Range("B15").Value = "Interested"
$ws.Range("C15").Value = "neutral"
$ws.Range("D15").Value = "2025-11-24 19:33:40"
$ws.Range("E15").Value = 919510038048

# Row 16 (new)
$ws.Range("A16").Value = "Yash"
$ws.Range("B16").Value = "Interested"
$ws.Range("C16").Value = "neutral"
$ws.Range("D16").Value = "2025-11-25 15:23:57"
$ws.Range("E16").Value = 917990747606

# Row 17 (new)
$ws.Range("A17").Value = "Batman"
$ws.Range("B17").Value = "Interested"
$ws.Range("C17").Value = "neutral"
$ws.Range("D17").Value = "2025-11-25 18:38:57"
$ws.Range("E17").Value = 917990747606

# Row 18 (new) - E18 keeps the leading "+" so it must stay text, not be
# coerced into a number by Excel's usual "+123" -> 123 auto-convert.
$ws.Range("A18").Value = "Yasha"
$ws.Range("B18").Value = "Interested"
$ws.Range("C18").Value = "neutral"
$ws.Range("D18").Value = "2025-11-25 19:31:00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "+917990747606"
$ws.Range("E18").Style = "Normal"
